$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C22").Value = 990
$ws.Range("D22").Value = 5732191
$ws.Range("E22").Value = 912.1882558879695
$ws.Range("G22").Value = 3.556485355648542
$ws.Range("H22").Value = 24.65729086021973
